$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data to the latest snapshot.
$ws.Range("D2").Value = "30.889.87"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "1.903.75"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.17"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4992"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2991"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06847"
$ws.Range("E9").Value = "  +3.75%  "
$ws.Range("D10").Value = "1.907.61"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.32"
$ws.Range("E11").Value = "  +3.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07337"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.99"
$ws.Range("E13").Value = "  +7.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.110"
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6832"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").Value = "30.885.37"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008074"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.43"
$ws.Range("E18").Value = "  +5.52%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "2.154.50"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.883"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "181.33"
$ws.Range("E23").Value = "  +34.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.102"
$ws.Range("E24").Value = "  +9.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.384"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.60"
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +11.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.953"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.397"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.382"
$ws.Range("E30").Value = "  +5.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08979"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.070"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05314"
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7517"
$ws.Range("E34").Value = "  +6.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.698"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01913"
$ws.Range("E37").Value = "  +16.35%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.198"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9409"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4395"
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.36"
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.884"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.773"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  +9.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05851"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.641"
$ws.Range("E48").Value = "  +5.27%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3930"
$ws.Range("E49").Value = "  +5.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.57"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.394"
$ws.Range("E51").Value = "  +4.16%  "
